$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing Doug Kelkhoff and subsequent rows down.
$ws.Rows.Item(9).Insert()

# Populate the new row with Jeff Thompson's info.
$ws.Range("A9").Value = "Jeff Thompson"
$ws.Range("B9").Value = "riskassessment"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "Workstream Co-lead"
$ws.Range("E9").Value = "person_placeholder.jpg"

# Update the active selection to match the saved view state.
$ws.Range("F9").Select() | Out-Null
